# Update overview.xlsx:
#  - Rework the "scalene" section on Sheet1 (rows 23-24/32 -> rows 24-28, plus new row 34)
#  - Add a new Sheet2 with the scalene-triangle JS snippet

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: clear out the old "scalene" rows (23, 24, 32) -----------------
# Row 32 ("scalene ", bold) has nothing below it, so deleting it first is a
# clean removal. Then row 23 ("or (This takes care..." , bold) is removed
# twice in a row so that both it and the row that shifts up into its place
# (old row 24, "(userInput1 === ...)") are gone, leaving rows 23-32 empty.
$ws1.Rows("32:32").Delete()
$ws1.Rows("23:23").Delete()
$ws1.Rows("23:23").Delete()

# --- Sheet1: write the new "scalene" rows -----------------------------------
$ws1.Range("A24").Value = "scalene "
$ws1.Range("A24").Font.Bold = $true

$ws1.Range("A25").Value = "For a scalene triangle, none of the sides is equal. However, the sum of any two sides of the triangle must be greater than the third side."
$ws1.Range("A26").Value = "userInput1!==userInput2 && userInput2!==userInput3 && userInput1!==userInput3"
$ws1.Range("A27").Value = "userInput1+userInput2 >userInput3 || userInput1+userInput3 >userInput2 || userInput3+userInput2 >userInput1"

# --- Sheet2: new worksheet with the scalene-triangle javascript snippet ----
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A2").Value = "// For a scalene triangle, none of the sides is equal. However, the sum of any two sides of the triangle must be greater than the third side."
$ws2.Range("A4").Value = 'var userInput1= parseInt(prompt("Please enter the 1st length of your triangle."));'
$ws2.Range("A5").Value = 'var userInput2= parseInt(prompt("Please enter the 2nd length of your triangle."));'
$ws2.Range("A6").Value = 'var userInput3= parseInt(prompt("Please enter the 3rd length of your triangle."));'
$ws2.Range("A7").Value = 'if((userInput1!==userInput2 && userInput2!==userInput3 && userInput1!==userInput3) && (userInput1+userInput2 >userInput3 || userInput1+userInput3 >userInput2 || userInput3+userInput2 >userInput1)&&(userInput1+userInput2 <=userInput3 || userInput1+userInput3 <=userInput2 || userInput3+userInput2 <=userInput1)){alert("Your triangle is a scalene triangle");'
$ws2.Range("A8").Value = "}"
$ws2.Range("A9").Value = 'else{alert("Your triangle is NOT a scalene triangle");}'
$ws2.Range("A11").Select()

# --- Sheet1: finish off the scalene section + add the new "4. Most..." row -
$ws1.Range("A28").Value = "userInput1+userInput2 >userInput3 && userInput1+userInput3 >userInput2 && userInput3+userInput2 >userInput1"

$ws1.Range("I7").Copy()
$ws1.Range("A34").PasteSpecial(-4122)
$ws1.Range("A34").Value = "4. Most importantly, if the sum of any two sides of the triangle is equal to or less than the third remaining side, then a triangle CANNOT be formed using those values. (For example, the values 9,4 & 3 cannot form a triangle.)"
$ws1.Rows("34:34").RowHeight = 60

# --- View state: restore Sheet1 as the active/selected tab -----------------
$ws1.Activate()
$ws1.Range("A25").Select()
